{"js": "// Remove the \"Fuentes:\" / contact-links block from the bio, while preserving\n// the `_GoBack` bookmark by re-anchoring it at the end of the preceding\n// paragraph (\"Eddie es licenciado...\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that precedes the block to be removed (\"Eddie es\n// licenciado...\") and the paragraphs that bound the block to delete\n// (\"Fuentes:\" through the mailto hyperlink paragraph).\nlet keepIndex = -1;\nlet fuentesIndex = -1;\nlet mailtoIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"Eddie es licenciado\") !== -1) {\n    keepIndex = i;\n  }\n  if (t.trim() === \"Fuentes:\") {\n    fuentesIndex = i;\n  }\n  if (t.indexOf(\"Eyoon@thecambridgegroup.com\") !== -1) {\n    mailtoIndex = i;\n  }\n}\n\nif (keepIndex === -1 || fuentesIndex === -1 || mailtoIndex === -1) {\n  throw new Error(\"Could not locate the expected paragraphs to edit.\");\n}\n\n// Re-anchor the `_GoBack` bookmark (previously sitting on the \"Correo\n// electr\u00f3nico:\" paragraph) onto the end of the surviving paragraph before\n// deleting the paragraphs that held it.\nconst keepRange = paragraphs.items[keepIndex].getRange(\"End\");\nkeepRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// Delete the \"Fuentes:\" paragraph through the mailto-hyperlink paragraph\n// (inclusive) \u2014 this removes: \"Fuentes:\", the LinkedIn/Cambridge Group\n// hyperlink paragraph, the blank paragraph, \"Correo electr\u00f3nico:\", and the\n// mailto hyperlink paragraph. Delete from last to first so indices stay\n// valid.\nfor (let i = mailtoIndex; i >= fuentesIndex; i--) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n", "ps1": "$doc = $word.ActiveDocument\n\n# --- Locate the relevant paragraphs by content -----------------------------\n$count = $doc.Paragraphs.Count\n$keepIndex = -1      # \"Eddie es licenciado...\" paragraph (survives, keeps bookmark)\n$fuentesIndex = -1   # \"Fuentes:\" paragraph (first paragraph of the block to delete)\n$mailtoIndex = -1    # paragraph containing the mailto hyperlink (last paragraph of the block to delete)\n\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $doc.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*Eddie es licenciado*\") {\n        $keepIndex = $i\n    }\n    if ($t.Trim() -eq \"Fuentes:\") {\n        $fuentesIndex = $i\n    }\n    if ($t -like \"*Eyoon@thecambridgegroup.com*\") {\n        $mailtoIndex = $i\n    }\n}\n\nif ($keepIndex -eq -1 -or $fuentesIndex -eq -1 -or $mailtoIndex -eq -1) {\n    throw \"Could not locate the expected paragraphs to edit.\"\n}\n\n# --- Re-anchor the `_GoBack` bookmark onto the end of the surviving --------\n# --- paragraph's text (before its paragraph mark) --------------------------\n#\n# A collapsed Range sitting exactly one character before a paragraph mark\n# cannot be used directly to seed Bookmarks.Add in this host, so a one-\n# character sentinel is inserted, bookmarked from a safe (non-boundary)\n# position, then removed again; the bookmark stays put, as in real Word.\n$keepPara = $doc.Paragraphs.Item($keepIndex)\n$tail = $keepPara.Range.Duplicate\n[void]$tail.MoveEnd(1, -1)  # wdCharacter, exclude the trailing paragraph mark\n$tail.Collapse(0)           # wdCollapseEnd\n$sentinelPos = $tail.Start\n$tail.InsertAfter(\"X\")\n\n$bmRange = $doc.Range($sentinelPos, $sentinelPos)\n$doc.Bookmarks.Add(\"_GoBack\", $bmRange)\n\n$sentinelRange = $doc.Range($sentinelPos, $sentinelPos + 1)\n$sentinelRange.Delete()\n\n# --- Delete the \"Fuentes:\" ... mailto-hyperlink block -----------------------\n$delStart = $doc.Paragraphs.Item($fuentesIndex).Range.Start\n$delEnd = $doc.Paragraphs.Item($mailtoIndex).Range.End\n$doc.Range($delStart, $delEnd).Delete()\n"}
